{"js": "// Update the worksheet date and the 25 two-digit multiplication problems\n// to the next day's generated set, matching the commit's regenerated\n// output (503736d).\nconst replacements = [\n  [\"2025-03-16 Sunday\", \"2025-03-17 Monday\"],\n  [\"37\u00d714=\", \"50\u00d744=\"],\n  [\"89\u00d752=\", \"70\u00d764=\"],\n  [\"71\u00d797=\", \"44\u00d733=\"],\n  [\"98\u00d752=\", \"50\u00d720=\"],\n  [\"50\u00d752=\", \"83\u00d730=\"],\n  [\"98\u00d793=\", \"39\u00d795=\"],\n  [\"25\u00d743=\", \"34\u00d713=\"],\n  [\"28\u00d720=\", \"18\u00d718=\"],\n  [\"98\u00d784=\", \"31\u00d718=\"],\n  [\"35\u00d772=\", \"12\u00d790=\"],\n  [\"48\u00d751=\", \"36\u00d736=\"],\n  [\"98\u00d795=\", \"29\u00d756=\"],\n  [\"52\u00d741=\", \"65\u00d796=\"],\n  [\"33\u00d745=\", \"48\u00d715=\"],\n  [\"98\u00d763=\", \"21\u00d786=\"],\n  [\"90\u00d763=\", \"49\u00d733=\"],\n  [\"94\u00d734=\", \"50\u00d794=\"],\n  [\"54\u00d799=\", \"88\u00d712=\"],\n  [\"53\u00d750=\", \"35\u00d760=\"],\n  [\"87\u00d767=\", \"72\u00d726=\"],\n  [\"78\u00d757=\", \"36\u00d750=\"],\n  [\"91\u00d754=\", \"20\u00d716=\"],\n  [\"65\u00d778=\", \"71\u00d764=\"],\n  [\"80\u00d751=\", \"69\u00d793=\"],\n  [\"76\u00d755=\", \"73\u00d783=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit multiplication problems\n# to the next day's generated set, matching the commit's regenerated\n# output (503736d).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-16 Sunday\", \"2025-03-17 Monday\"),\n    @(\"37\u00d714=\", \"50\u00d744=\"),\n    @(\"89\u00d752=\", \"70\u00d764=\"),\n    @(\"71\u00d797=\", \"44\u00d733=\"),\n    @(\"98\u00d752=\", \"50\u00d720=\"),\n    @(\"50\u00d752=\", \"83\u00d730=\"),\n    @(\"98\u00d793=\", \"39\u00d795=\"),\n    @(\"25\u00d743=\", \"34\u00d713=\"),\n    @(\"28\u00d720=\", \"18\u00d718=\"),\n    @(\"98\u00d784=\", \"31\u00d718=\"),\n    @(\"35\u00d772=\", \"12\u00d790=\"),\n    @(\"48\u00d751=\", \"36\u00d736=\"),\n    @(\"98\u00d795=\", \"29\u00d756=\"),\n    @(\"52\u00d741=\", \"65\u00d796=\"),\n    @(\"33\u00d745=\", \"48\u00d715=\"),\n    @(\"98\u00d763=\", \"21\u00d786=\"),\n    @(\"90\u00d763=\", \"49\u00d733=\"),\n    @(\"94\u00d734=\", \"50\u00d794=\"),\n    @(\"54\u00d799=\", \"88\u00d712=\"),\n    @(\"53\u00d750=\", \"35\u00d760=\"),\n    @(\"87\u00d767=\", \"72\u00d726=\"),\n    @(\"78\u00d757=\", \"36\u00d750=\"),\n    @(\"91\u00d754=\", \"20\u00d716=\"),\n    @(\"65\u00d778=\", \"71\u00d764=\"),\n    @(\"80\u00d751=\", \"69\u00d793=\"),\n    @(\"76\u00d755=\", \"73\u00d783=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
